$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 733.1667
$ws.Range("I2").Value = 739.8
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 739.8
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -626.8
$ws.Range("N2").Value = -926
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H69").Value = 12500
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 12500
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H96").Value = 681.5714
$ws.Range("I96").Value = 388
$ws.Range("J96").Value = 799
$ws.Range("K96").Value = 1164
$ws.Range("L96").Value = 2397
$ws.Range("M96").Value = 209
$ws.Range("N96").Value = -5143
$ws.Range("H116").Value = 4867.3335
$ws.Range("I116").Value = 4734.6665
$ws.Range("K116").Value = 4734.6665
$ws.Range("M116").Value = -1292.6665
$ws.Range("H125").Value = 810.6667
$ws.Range("I125").Value = 766
$ws.Range("J125").Value = 900
$ws.Range("K125").Value = 6894
$ws.Range("L125").Value = 8100
$ws.Range("M125").Value = -4434
$ws.Range("N125").Value = -13020

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2639.5
$ws.Range("I45").Value = 2639.5
$ws.Range("K45").Value = 2639.5
$ws.Range("M45").Value = -2262.5
$ws.Range("H119").Value = 66527.14
$ws.Range("J119").Value = 66527.14
$ws.Range("L119").Value = 66527.14
$ws.Range("N119").Value = -76203.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1435.3334
$ws.Range("I20").Value = 1053
$ws.Range("J20").Value = 2200
$ws.Range("K20").Value = 1053
$ws.Range("L20").Value = 2200
$ws.Range("M20").Value = -806
$ws.Range("N20").Value = -2694
$ws.Range("H29").Value = 1216.619
$ws.Range("I29").Value = 1329.8
$ws.Range("J29").Value = 1181.25
$ws.Range("K29").Value = 1329.8
$ws.Range("L29").Value = 1181.25
$ws.Range("M29").Value = -1040.8
$ws.Range("N29").Value = -1759.25
$ws.Range("H37").Value = 426
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H134").Value = 706.8
$ws.Range("I134").Value = 706.8
$ws.Range("K134").Value = 2120.4
$ws.Range("M134").Value = 414.6000000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 195.85715
$ws.Range("I7").Value = 221.88889
$ws.Range("K7").Value = 221.88889
$ws.Range("M7").Value = -108.88889
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 706.8182
$ws.Range("I22").Value = 810.1111
$ws.Range("J22").Value = 242
$ws.Range("K22").Value = 810.1111
$ws.Range("L22").Value = 242
$ws.Range("M22").Value = -460.1111
$ws.Range("N22").Value = -942
$ws.Range("H86").Value = 200012000
$ws.Range("I86").Value = 333341660
$ws.Range("J86").Value = 17499
$ws.Range("K86").Value = 333341660
$ws.Range("L86").Value = 17499
$ws.Range("M86").Value = -333340537
$ws.Range("N86").Value = -19745
$ws.Range("H89").Value = 200012000
$ws.Range("I89").Value = 333341660
$ws.Range("J89").Value = 17499
$ws.Range("K89").Value = 1666708300
$ws.Range("L89").Value = 87495
$ws.Range("M89").Value = -1666702684
$ws.Range("N89").Value = -98727
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2139
$ws.Range("I132").Value = 1835.2222
$ws.Range("K132").Value = 5505.6666
$ws.Range("M132").Value = -2975.6666
$ws.Range("H134").Value = 1546.7142
$ws.Range("J134").Value = 1765.6
$ws.Range("L134").Value = 5296.799999999999
$ws.Range("N134").Value = -10366.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H33").Value = 597.7143
$ws.Range("I33").Value = 617
$ws.Range("K33").Value = 3702
$ws.Range("M33").Value = -3419
$ws.Range("H80").Value = 5958.3335
$ws.Range("I80").Value = 5937.5
$ws.Range("K80").Value = 17812.5
$ws.Range("M80").Value = -16876.5
$ws.Range("H83").Value = 5958.3335
$ws.Range("I83").Value = 5937.5
$ws.Range("K83").Value = 53437.5
$ws.Range("M83").Value = -48757.5
$ws.Range("H114").Value = 1559.8
$ws.Range("J114").Value = 1433
$ws.Range("L114").Value = 4299
$ws.Range("N114").Value = -10807

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 32000
$ws.Range("J63").Value = 32000
$ws.Range("L63").Value = 32000
$ws.Range("N63").Value = -33372
$ws.Range("H66").Value = 32000
$ws.Range("J66").Value = 32000
$ws.Range("L66").Value = 96000
$ws.Range("N66").Value = -102864
$ws.Range("H102").Value = 687.2727
$ws.Range("I102").Value = 729
$ws.Range("J102").Value = 499.5
$ws.Range("K102").Value = 729
$ws.Range("L102").Value = 499.5
$ws.Range("M102").Value = 893
$ws.Range("N102").Value = -3743.5
$ws.Range("H132").Value = 1854.4286
$ws.Range("I132").Value = 1821.25
$ws.Range("J132").Value = 1898.6666
$ws.Range("K132").Value = 5463.75
$ws.Range("L132").Value = 5695.9998
$ws.Range("M132").Value = -2933.75
$ws.Range("N132").Value = -10755.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8770.4
$ws.Range("I7").Value = 9225
$ws.Range("J7").Value = 8467.333000000001
$ws.Range("K7").Value = 9225
$ws.Range("L7").Value = 8467.333000000001
$ws.Range("M7").Value = -9113
$ws.Range("N7").Value = -8691.333000000001
$ws.Range("H61").Value = 3000
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3404
$ws.Range("H103").Value = 15549.25
$ws.Range("J103").Value = 15549.25
$ws.Range("L103").Value = 15549.25
$ws.Range("N103").Value = -17893.25
$ws.Range("H113").Value = 3000
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 8770.4
$ws.Range("I126").Value = 9225
$ws.Range("J126").Value = 8467.333000000001
$ws.Range("K126").Value = 27675
$ws.Range("L126").Value = 25401.999
$ws.Range("M126").Value = -25205
$ws.Range("N126").Value = -30341.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 255.45454
$ws.Range("I113").Value = 170.5
$ws.Range("K113").Value = 511.5
$ws.Range("M113").Value = 1658.5
$ws.Range("H117").Value = 32000
$ws.Range("J117").Value = 32000
$ws.Range("L117").Value = 32000
$ws.Range("N117").Value = -41178
$ws.Range("H132").Value = 2585.625
$ws.Range("I132").Value = 2377
$ws.Range("J132").Value = 2933.3333
$ws.Range("K132").Value = 7131
$ws.Range("L132").Value = 8799.999899999999
$ws.Range("M132").Value = -4601
$ws.Range("N132").Value = -13859.9999
